$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F2").Value = 5360
$ws.Range("F3").Value = 20500
$ws.Range("F4").Value = 892
$ws.Range("F5").Value = 1544

$ws.Range("F6").Select()
